$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$ws.Range("A1").Value2 = 45309
$ws.Range("D32").Value2 = 2033
$ws.Range("D33").Value2 = 2430
$ws.Range("D34").Value2 = 2720
$ws.Range("D35").Value2 = 2930
$ws.Range("D36").Value2 = 3100
$ws.Range("D37").Value2 = 3481
$ws.Range("D38").Value2 = 3830
$ws.Range("D39").Value2 = 4100
$ws.Range("D40").Value2 = 4320
$ws.Range("D46").Value2 = 2480
$ws.Range("D47").Value2 = 2770
$ws.Range("D48").Value2 = 3040
$ws.Range("D49").Value2 = 3540
$ws.Range("D50").Value2 = 3851
$ws.Range("D51").Value2 = 4290
$ws.Range("D52").Value2 = 4816
$ws.Range("D53").Value2 = 5020
$ws.Range("D54").Value2 = 5540
$ws.Range("D55").Value2 = 6190
$ws.Range("D56").Value2 = 6721
$ws.Range("D57").Value2 = 7550
$ws.Range("D58").Value2 = 8260
$ws.Range("D59").Value2 = 9050
$ws.Range("D60").Value2 = 10317.112
$ws.Range("D61").Value2 = 10949.764
$ws.Range("D67").Value2 = 4590
$ws.Range("D68").Value2 = 5235
$ws.Range("D69").Value2 = 5640
$ws.Range("D70").Value2 = 5730
$ws.Range("D71").Value2 = 6470
$ws.Range("D72").Value2 = 6970
$ws.Range("D73").Value2 = 7762
$ws.Range("D74").Value2 = 8335
$ws.Range("D75").Value2 = 9675
$ws.Range("D76").Value2 = 10520
$ws.Range("D77").Value2 = 11710
$ws.Range("D78").Value2 = 12890
$ws.Range("D79").Value2 = 14300
$ws.Range("D80").Value2 = 19730
$ws.Range("D81").Value2 = 21980
$ws.Range("D87").Value2 = 8140
$ws.Range("D88").Value2 = 8400
$ws.Range("D89").Value2 = 9840
$ws.Range("D90").Value2 = 10780
$ws.Range("D91").Value2 = 11900
$ws.Range("D92").Value2 = 14540
$ws.Range("D93").Value2 = 15130
$ws.Range("D94").Value2 = 16800
$ws.Range("D95").Value2 = 18400
$ws.Range("D96").Value2 = 20350
